# Update "想去人数" (want-to-go count) figures in the 展览 (Exhibition),
# 演出 (Performance) and 全部类型 (All types) sheets, matching the data
# refresh captured in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (sheetId 1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 5463
$wsExhibit.Range("F6").Value = 824
$wsExhibit.Range("F7").Value = 24
$wsExhibit.Range("F8").Value = 350

# 演出 (sheetId 2)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 44

# 全部类型 (sheetId 4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 5464
$wsAll.Range("F6").Value = 824
$wsAll.Range("F7").Value = 24
$wsAll.Range("F8").Value = 44
$wsAll.Range("F9").Value = 350
